# Auto-generated edit script: updates computed profit columns (H-N)
# for specific Leve rows across multiple sheets, per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(7, 8).Value = 1650
$ws.Cells.Item(7, 9).Value = 366.66666
$ws.Cells.Item(7, 10).Value = 2933.3333
$ws.Cells.Item(7, 11).Value = 366.66666
$ws.Cells.Item(7, 12).Value = 2933.3333
$ws.Cells.Item(7, 13).Value = -254.66666
$ws.Cells.Item(7, 14).Value = -3157.3333

$ws.Cells.Item(14, 8).Value = 1650
$ws.Cells.Item(14, 9).Value = 366.66666
$ws.Cells.Item(14, 10).Value = 2933.3333
$ws.Cells.Item(14, 11).Value = 366.66666
$ws.Cells.Item(14, 12).Value = 2933.3333
$ws.Cells.Item(14, 13).Value = -175.66666
$ws.Cells.Item(14, 14).Value = -3315.3333

$ws.Cells.Item(62, 8).Value = 4665.6665
$ws.Cells.Item(62, 9).Value = 1998.5
$ws.Cells.Item(62, 10).Value = 10000
$ws.Cells.Item(62, 11).Value = 1998.5
$ws.Cells.Item(62, 12).Value = 10000
$ws.Cells.Item(62, 13).Value = -1374.5
$ws.Cells.Item(62, 14).Value = -11248

$ws.Cells.Item(65, 8).Value = 4665.6665
$ws.Cells.Item(65, 9).Value = 1998.5
$ws.Cells.Item(65, 10).Value = 10000
$ws.Cells.Item(65, 11).Value = 9992.5
$ws.Cells.Item(65, 12).Value = 50000
$ws.Cells.Item(65, 13).Value = -6872.5
$ws.Cells.Item(65, 14).Value = -56240

$ws.Cells.Item(101, 8).Value = 3378
$ws.Cells.Item(101, 9).Value = 4142
$ws.Cells.Item(101, 10).Value = 1850
$ws.Cells.Item(101, 11).Value = 12426
$ws.Cells.Item(101, 12).Value = 5550
$ws.Cells.Item(101, 13).Value = -10804
$ws.Cells.Item(101, 14).Value = -8794

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1989.4
$ws.Cells.Item(2, 9).Value = 1989.4
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 1989.4
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -1876.4
$ws.Cells.Item(2, 14).ClearContents()

$ws.Cells.Item(15, 8).Value = 15000
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 15000
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 15000
$ws.Cells.Item(15, 14).Value = -15700

$ws.Cells.Item(102, 8).Value = 558.1818
$ws.Cells.Item(102, 9).Value = 589
$ws.Cells.Item(102, 10).Value = 250
$ws.Cells.Item(102, 11).Value = 589
$ws.Cells.Item(102, 12).Value = 250
$ws.Cells.Item(102, 13).Value = 1033
$ws.Cells.Item(102, 14).Value = -3494

$ws.Cells.Item(116, 8).Value = 1989.4
$ws.Cells.Item(116, 9).Value = 1989.4
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 1989.4
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = 304.5999999999999
$ws.Cells.Item(116, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1989.4
$ws.Cells.Item(3, 9).Value = 1989.4
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 1989.4
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = -1875.4
$ws.Cells.Item(3, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 10075212
$ws.Cells.Item(6, 9).Value = 10075212
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 10075212
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -10075099

$ws.Cells.Item(10, 8).Value = 120.09091
$ws.Cells.Item(10, 9).Value = 121.3
$ws.Cells.Item(10, 10).Value = 108
$ws.Cells.Item(10, 11).Value = 121.3
$ws.Cells.Item(10, 12).Value = 108
$ws.Cells.Item(10, 13).Value = 17.7
$ws.Cells.Item(10, 14).Value = -386

$ws.Cells.Item(22, 8).Value = 503
$ws.Cells.Item(22, 9).Value = 465.15384
$ws.Cells.Item(22, 10).Value = 667
$ws.Cells.Item(22, 11).Value = 465.15384
$ws.Cells.Item(22, 12).Value = 667
$ws.Cells.Item(22, 13).Value = -115.15384
$ws.Cells.Item(22, 14).Value = -1367

$ws.Cells.Item(31, 8).Value = 4980
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 4980
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 4980
$ws.Cells.Item(31, 13).ClearContents()
$ws.Cells.Item(31, 14).Value = -5570

$ws.Cells.Item(34, 8).Value = 4980
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 4980
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 4980
$ws.Cells.Item(34, 13).ClearContents()
$ws.Cells.Item(34, 14).Value = -5384

$ws.Cells.Item(58, 8).Value = 3670.3333
$ws.Cells.Item(58, 9).Value = 1012
$ws.Cells.Item(58, 10).Value = 4999.5
$ws.Cells.Item(58, 11).Value = 1012
$ws.Cells.Item(58, 12).Value = 4999.5
$ws.Cells.Item(58, 13).Value = -809
$ws.Cells.Item(58, 14).Value = -5405.5

$ws.Cells.Item(131, 8).Value = 45666
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 45666
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 45666
$ws.Cells.Item(131, 14).Value = -55746

$ws.Cells.Item(136, 8).Value = 3670.3333
$ws.Cells.Item(136, 9).Value = 1012
$ws.Cells.Item(136, 10).Value = 4999.5
$ws.Cells.Item(136, 11).Value = 3036
$ws.Cells.Item(136, 12).Value = 14998.5
$ws.Cells.Item(136, 13).Value = -486
$ws.Cells.Item(136, 14).Value = -20098.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 12.25
$ws.Cells.Item(7, 9).Value = 10
$ws.Cells.Item(7, 10).Value = 19
$ws.Cells.Item(7, 11).Value = 30
$ws.Cells.Item(7, 12).Value = 57
$ws.Cells.Item(7, 13).Value = 82
$ws.Cells.Item(7, 14).Value = -281

$ws.Cells.Item(107, 8).Value = 747.25
$ws.Cells.Item(107, 9).Value = 497.5
$ws.Cells.Item(107, 10).Value = 997
$ws.Cells.Item(107, 11).Value = 1492.5
$ws.Cells.Item(107, 12).Value = 2991
$ws.Cells.Item(107, 13).Value = 427.5
$ws.Cells.Item(107, 14).Value = -6831

$ws.Cells.Item(141, 8).Value = 3000
$ws.Cells.Item(141, 9).Value = 2600
$ws.Cells.Item(141, 10).Value = 3200
$ws.Cells.Item(141, 11).Value = 7800
$ws.Cells.Item(141, 12).Value = 9600
$ws.Cells.Item(141, 13).Value = -2620
$ws.Cells.Item(141, 14).Value = -19960

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(17, 8).Value = 450
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 450
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 450
$ws.Cells.Item(17, 13).ClearContents()
$ws.Cells.Item(17, 14).Value = -786

$ws.Cells.Item(80, 8).Value = 12639
$ws.Cells.Item(80, 9).Value = 3323.75
$ws.Cells.Item(80, 10).Value = 49900
$ws.Cells.Item(80, 11).Value = 3323.75
$ws.Cells.Item(80, 12).Value = 49900
$ws.Cells.Item(80, 13).Value = -2325.75
$ws.Cells.Item(80, 14).Value = -51896

$ws.Cells.Item(83, 8).Value = 12639
$ws.Cells.Item(83, 9).Value = 3323.75
$ws.Cells.Item(83, 10).Value = 49900
$ws.Cells.Item(83, 11).Value = 16618.75
$ws.Cells.Item(83, 12).Value = 249500
$ws.Cells.Item(83, 13).Value = -11626.75
$ws.Cells.Item(83, 14).Value = -259484

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 5411
$ws.Cells.Item(9, 9).Value = 133
$ws.Cells.Item(9, 10).Value = 8050
$ws.Cells.Item(9, 11).Value = 133
$ws.Cells.Item(9, 12).Value = 8050
$ws.Cells.Item(9, 13).Value = 91
$ws.Cells.Item(9, 14).Value = -8498

$ws.Cells.Item(11, 8).Value = 1499
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 1499
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 1499
$ws.Cells.Item(11, 13).ClearContents()
$ws.Cells.Item(11, 14).Value = -1779

$ws.Cells.Item(13, 8).Value = 19000
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 19000
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 19000
$ws.Cells.Item(13, 13).ClearContents()
$ws.Cells.Item(13, 14).Value = -19280

$ws.Cells.Item(17, 8).Value = 3499
$ws.Cells.Item(17, 9).Value = 2497
$ws.Cells.Item(17, 10).Value = 4000
$ws.Cells.Item(17, 11).Value = 2497
$ws.Cells.Item(17, 12).Value = 4000
$ws.Cells.Item(17, 13).Value = -2327
$ws.Cells.Item(17, 14).Value = -4340

$ws.Cells.Item(82, 8).Value = 3566.6667
$ws.Cells.Item(82, 9).Value = 2700
$ws.Cells.Item(82, 10).Value = 4000
$ws.Cells.Item(82, 11).Value = 2700
$ws.Cells.Item(82, 12).Value = 4000
$ws.Cells.Item(82, 13).Value = -2339
$ws.Cells.Item(82, 14).Value = -4722

$ws.Cells.Item(85, 8).Value = 3566.6667
$ws.Cells.Item(85, 9).Value = 2700
$ws.Cells.Item(85, 10).Value = 4000
$ws.Cells.Item(85, 11).Value = 2700
$ws.Cells.Item(85, 12).Value = 4000
$ws.Cells.Item(85, 13).Value = -1452
$ws.Cells.Item(85, 14).Value = -6496

Write-Output "Applied profit-column updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets"